# The NATMI pipeline was re-run with updated TPM values for the
# Cxcl10 -> Ccr3 ligand/receptor pair (YoungD2, lrc2p). This refreshes the
# derived ligand/receptor/edge expression & specificity columns (G:J, M:T,
# plus the few K/L receptor-expressing-cell counts that shift because the
# underlying detection threshold moved) on rows 2-19 to the newly computed
# values, while leaving the categorical columns (A:F) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ "G" = 8.850858000000001; "H" = 26.552574; "I" = 0.05442939716240135; "J" = 0.05442939716240137; "M" = 0.1808983333333333; "N" = 0.5426949999999999; "O" = 0.09546831801815302; "P" = 0.09546831801815302; "Q" = 1.60110546077; "R" = 14.40994914693; "S" = 0.005196282997836488; "T" = 0.005196282997836489 }
    3 = @{ "G" = 8.850858000000001; "H" = 26.552574; "I" = 0.05442939716240135; "J" = 0.05442939716240137; "M" = 1.572737; "N" = 4.718211; "O" = 0.8300051930177132; "P" = 0.8300051930177132; "Q" = 13.920071858346; "R" = 125.280646725114; "S" = 0.04517668229761671; "T" = 0.04517668229761672 }
    4 = @{ "G" = 8.850858000000001; "H" = 26.552574; "I" = 0.05442939716240135; "J" = 0.05442939716240137; "K" = 1; "L" = 0.3333333333333333; "M" = 0.1412166666666667; "N" = 0.42365; "O" = 0.07452648896413371; "P" = 0.07452648896413369; "Q" = 1.2498886639; "R" = 11.2489979751; "S" = 0.004056431866948155; "T" = 0.004056431866948155 }
    5 = @{ "I" = 0.2412735821509021; "J" = 0.2412735821509022; "M" = 0.1808983333333333; "N" = 0.5426949999999999; "O" = 0.09546831801815302; "P" = 0.09546831801815302; "Q" = 7.097349411545554; "R" = 63.87614470390999; "S" = 0.02303398307016129; "T" = 0.0230339830701613 }
    6 = @{ "I" = 0.2412735821509021; "J" = 0.2412735821509022; "M" = 1.572737; "N" = 4.718211; "O" = 0.8300051930177132; "P" = 0.8300051930177132; "Q" = 61.70462610563534; "R" = 555.341634950718; "S" = 0.2002583261232346; "T" = 0.2002583261232347 }
    7 = @{ "I" = 0.2412735821509021; "J" = 0.2412735821509022; "K" = 1; "L" = 0.3333333333333333; "M" = 0.1412166666666667; "N" = 0.42365; "O" = 0.07452648896413371; "P" = 0.07452648896413369; "Q" = 5.540482367077779; "R" = 49.86434130370001; "S" = 0.01798127295750621; "T" = 0.01798127295750622 }
    8 = @{ "G" = 50.430027; "H" = 151.290081; "I" = 0.3101254102702387; "J" = 0.3101254102702387; "M" = 0.1808983333333333; "N" = 0.5426949999999999; "O" = 0.09546831801815302; "P" = 0.09546831801815302; "Q" = 9.122707834254999; "R" = 82.10437050829499; "S" = 0.02960715129318932; "T" = 0.02960715129318933 }
    9 = @{ "G" = 50.430027; "H" = 151.290081; "I" = 0.3101254102702387; "J" = 0.3101254102702387; "M" = 1.572737; "N" = 4.718211; "O" = 0.8300051930177132; "P" = 0.8300051930177132; "Q" = 79.313169373899; "R" = 713.8185243650911; "S" = 0.2574057010110469; "T" = 0.257405701011047 }
    10 = @{ "G" = 50.430027; "H" = 151.290081; "I" = 0.3101254102702387; "J" = 0.3101254102702387; "K" = 1; "L" = 0.3333333333333333; "M" = 0.1412166666666667; "N" = 0.42365; "O" = 0.07452648896413371; "P" = 0.07452648896413369; "Q" = 7.121560312850002; "R" = 64.09404281565001; "S" = 0.02311255796600238; "T" = 0.02311255796600238 }
    11 = @{ "G" = 3.269985333333333; "H" = 9.809956; "I" = 0.02010916121614733; "J" = 0.02010916121614734; "M" = 0.1808983333333333; "N" = 0.5426949999999999; "O" = 0.09546831801815302; "P" = 0.09546831801815302; "Q" = 0.5915348968244443; "R" = 5.323814071419999; "S" = 0.001919787798061462; "T" = 0.001919787798061463 }
    12 = @{ "G" = 3.269985333333333; "H" = 9.809956; "I" = 0.02010916121614733; "J" = 0.02010916121614734; "M" = 1.572737; "N" = 4.718211; "O" = 0.8300051930177132; "P" = 0.8300051930177132; "Q" = 5.142826923190666; "R" = 46.285442308716; "S" = 0.01669070823663268; "T" = 0.01669070823663269 }
    13 = @{ "G" = 3.269985333333333; "H" = 9.809956; "I" = 0.02010916121614733; "J" = 0.02010916121614734; "K" = 1; "L" = 0.3333333333333333; "M" = 0.1412166666666667; "N" = 0.42365; "O" = 0.07452648896413371; "P" = 0.07452648896413369; "Q" = 0.4617764288222223; "R" = 4.1559878594; "S" = 0.00149866518145319; "T" = 0.00149866518145319 }
    14 = @{ "G" = 47.074941; "H" = 141.224823; "I" = 0.2894929124482182; "J" = 0.2894929124482182; "M" = 0.1808983333333333; "N" = 0.5426949999999999; "O" = 0.09546831801815302; "P" = 0.09546831801815302; "Q" = 8.515778368664998; "R" = 76.642005317985; "S" = 0.02763740142960782; "T" = 0.02763740142960783 }
    15 = @{ "G" = 47.074941; "H" = 141.224823; "I" = 0.2894929124482182; "J" = 0.2894929124482182; "M" = 1.572737; "N" = 4.718211; "O" = 0.8300051930177132; "P" = 0.8300051930177132; "Q" = 74.036501483517; "R" = 666.3285133516531; "S" = 0.2402806206738433; "T" = 0.2402806206738433 }
    16 = @{ "G" = 47.074941; "H" = 141.224823; "I" = 0.2894929124482182; "J" = 0.2894929124482182; "K" = 1; "L" = 0.3333333333333333; "M" = 0.1412166666666667; "N" = 0.42365; "O" = 0.07452648896413371; "P" = 0.07452648896413369; "Q" = 6.647766251550001; "R" = 59.82989626395001; "S" = 0.02157489034476706; "T" = 0.02157489034476706 }
    17 = @{ "G" = 13.751998; "H" = 41.255994; "I" = 0.08456953675209218; "J" = 0.08456953675209219; "M" = 0.1808983333333333; "N" = 0.5426949999999999; "O" = 0.09546831801815302; "P" = 0.09546831801815302; "Q" = 2.487713518203333; "R" = 22.38942166383; "S" = 0.008073711429296616; "T" = 0.008073711429296616 }
    18 = @{ "G" = 13.751998; "H" = 41.255994; "I" = 0.08456953675209218; "J" = 0.08456953675209219; "M" = 1.572737; "N" = 4.718211; "O" = 0.8300051930177132; "P" = 0.8300051930177132; "Q" = 21.628276078526; "R" = 194.654484706734; "S" = 0.07019315467533886; "T" = 0.07019315467533888 }
    19 = @{ "G" = 13.751998; "H" = 41.255994; "I" = 0.08456953675209218; "J" = 0.08456953675209219; "K" = 1; "L" = 0.3333333333333333; "M" = 0.1412166666666667; "N" = 0.42365; "O" = 0.07452648896413371; "P" = 0.07452648896413369; "Q" = 1.942011317566667; "R" = 17.4781018581; "S" = 0.006302670647456697; "T" = 0.006302670647456697 }
}
foreach ($rowKey in $newValues.Keys) {
    $rowValues = $newValues[$rowKey]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$rowKey").Value = $rowValues[$col]
    }
}
